# Remove the "Pick a date range" / "date-range" question row from the
# questions sheet (row 8). Excel shifts row 9 ("Pick a time (24 hrs)")
# up into row 8 and automatically prunes the now-unused shared strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(8).Delete()

# Move the active selection, matching the saved view state.
$ws.Range("D15").Select()
